# Update cryptocurrency price/volume data on Sheet1 to reflect the latest
# scrape performed by the GitHub Actions workflow.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new value. All of these columns are stored as
# plain text (inline strings) in the workbook, so we force the values in
# as text (prefixed with a leading apostrophe-equivalent via NumberFormat)
# to avoid Excel silently reinterpreting percentages/decimals as numbers
# and altering their displayed precision (e.g. dropping trailing zeros).

$updates = @{
    "D2"  = "307.37";  "E2"  = "1.71%"
    "D3"  = "36.27";   "E3"  = "3.19%"
    "D4"  = "5.097";   "E4"  = "0.85%"
    "D5"  = "0.08127"; "E5"  = "2.83%"
    "D6"  = "1.946";   "E6"  = "0.22%"
    "D7"  = "7.781";   "E7"  = "0.38%"
    "D8"  = "0.9308";  "E8"  = "0.81%"
    "D9"  = "0.1398";  "E9"  = "16.05%"
    "D10" = "0.1926";  "E10" = "4.81%"
    "D11" = "0.09241"; "E11" = "-1.69%"
    "D12" = "0.03524"; "E12" = "-0.55%"
    "D13" = "0.09866"; "E13" = "-0.04%"
    "D14" = "0.001410";"E14" = "1.01%"
    "D15" = "0.005757";"E15" = "-1.63%"
    "D16" = "3.605";   "E16" = "3.23%"
    "D17" = "4.185";   "E17" = "3.99%"
    "E18" = "4.03%"
    "D19" = "0.3438"
    "E20" = "1.02%"
    "D21" = "4.897";   "E21" = "-2.76%"
    "D22" = "0.2499";  "E22" = "4.08%"
    "D23" = "0.04507"; "E23" = "-0.06%"
    "D24" = "0.001217";"E24" = "0.25%"
    "E25" = "6.54%"
    "E26" = "-0.77%"
    "D39" = "0.02010"; "E39" = "6.04%"
    "D40" = "0.04941"; "E40" = "4.93%"
    "D41" = "0.007662";"E41" = "0.39%"
    "D42" = "0.01022"; "E42" = "7.05%"
    "E43" = "4.42%"
    "D44" = "0.002102";"E44" = "-0.45%"
    "D45" = "0.01032"; "E45" = "-7.86%"
    "D46" = "0.00006428"; "E46" = "6.73%"
    "D47" = "0.00000000751"; "E47" = "0.03%"
    "E50" = "0.03%"
    "E51" = "0.03%"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
